$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A14").Value = 13
$ws.Range("B14").Value = 62.5

$ws.Range("A15").Select()
